$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear out the contents of B5 and D5 (the "app.fancy-slider.controls" row),
# marking the module as integrated, and restore the default "Good" (done) style
# to match the other completed rows.
$ws.Range("B5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("B5").Style = "Good"
$ws.Range("D5").Style = "Good"

# Move the active selection to B9, as recorded in the saved view state.
$ws.Range("B9").Select()
